# Enhance Azure DevOps integration and improve chatbot functionality
#
# Applies the following changes to the "work_items_due_dates" workbook:
#  - Sheet "This Friday": remove the "TEST_03: Coffee Cake" work item, reorder the
#    remaining TEST_ items, flip their State to "In Progress", and re-assign
#    "TEST_05: Strawberry" to Priththiha Nemikumar.
#  - Sheet "Next Friday": remove the "QA | Sure, whatever" work item.
#  - Sheet "Friday After Next": move "Bug | Minor | Performance Degradation" to
#    "In Progress" and replace the "Arian Fooladray" assignee with his Azure DevOps
#    login e-mail (afooladray@fgfbrands.com).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "This Friday"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("This Friday")

# Drop the "TEST_03: Coffee Cake" row (row 4: ID 986) entirely - shifts rows up.
$ws1.Rows.Item(4).Delete()

# Rewrite the remaining data rows (2-9) in their final order/content.
$sheet1Data = @(
    @(985, "TEST_02: Something", "In Progress", "Kenneth Kui", "2024-12-03T00:00:00Z", "Preet Patel", 2),
    @(984, "TEST_01: Introduction", "In Progress", "Kenneth Kui", "2024-12-02T00:00:00Z", "Preet Patel", 2),
    @(988, "TEST_05: Strawberry", "In Progress", "Priththiha Nemikumar", "2024-12-06T00:00:00Z", "Preet Patel", 2),
    @(987, "TEST_04: Testing Sandwiches", "In Progress", "Kenneth Kui", "2024-12-05T00:00:00Z", "Preet Patel", 2),
    @(992, "Creating BRD", "New", "Kamini Patel", "2024-12-06T00:00:00Z", "Preet Patel", 2),
    @(990, "Planning", "In Progress", "Kenneth Kui", "2024-12-06T00:00:00Z", "Preet Patel", 2),
    @(1005, "Design Start", "New", "Kamini Patel", "2024-12-02T00:00:00Z", "Preet Patel", 2),
    @(991, "Gather Design", "In Progress", "Kenneth Kui", "2024-12-06T00:00:00Z", "Preet Patel", 2)
)

$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: "Next Friday"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Next Friday")

# Drop the "QA | Sure, whatever" row (row 2: ID 945) entirely - shifts rows up.
$ws2.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# Sheet 3: "Friday After Next"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Friday After Next")

# "Bug | Minor | Performance Degradation" moves to In Progress.
$ws3.Cells.Item(2, 3).Value = "In Progress"

# Replace "Arian Fooladray" with his e-mail address across the sheet.
for ($r = 2; $r -le 11; $r++) {
    if ($ws3.Cells.Item($r, 4).Value2 -eq "Arian Fooladray") {
        $ws3.Cells.Item($r, 4).Value = "afooladray@fgfbrands.com"
    }
}
